# Coinranking cryptos list refresh -- Wed Apr 19 11:39:58 UTC 2023 (GitHub Actions)
# Updates the "Price" (D) / "Volume(1h)" (E) columns for rows 2-51 to the
# scraper's latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "29.482.87"
$ws.Range('E2').Value = "  -2.70%  "

$ws.Range('D3').Value = "1.994.62"
$ws.Range('E3').Value = "  -6.14%  "

$ws.Range('D4').Value = "'1.006"
$ws.Range('E4').Value = "  +0.04%  "

$ws.Range('D5').Value = "'329.67"
$ws.Range('E5').Value = "  -5.09%  "

$ws.Range('D6').Value = "'1.006"
$ws.Range('E6').Value = "  +0.10%  "

$ws.Range('D7').Value = "'0.5016"
$ws.Range('E7').Value = "  -3.95%  "

$ws.Range('D8').Value = "'0.4226"
$ws.Range('E8').Value = "  -5.72%  "

$ws.Range('D9').Value = "'52.02"
$ws.Range('E9').Value = "  -4.03%  "

$ws.Range('D10').Value = "'0.08918"
$ws.Range('E10').Value = "  -5.03%  "

$ws.Range('D11').Value = "'1.123"

$ws.Range('D12').Value = "'23.40"
$ws.Range('E12').Value = "  -8.21%  "

$ws.Range('D13').Value = "'8.105"
$ws.Range('E13').Value = "  -6.78%  "

$ws.Range('D14').Value = "1.995.63"
$ws.Range('E14').Value = "  -5.95%  "

$ws.Range('D15').Value = "'6.521"
$ws.Range('E15').Value = "  -6.56%  "

$ws.Range('D16').Value = "'96.32"
$ws.Range('E16').Value = "  -6.50%  "

$ws.Range('D17').Value = "'1.007"
$ws.Range('E17').Value = "  +0.06%  "

$ws.Range('E18').Value = "  -6.07%  "

$ws.Range('D19').Value = "'0.06626"
$ws.Range('E19').Value = "  -1.18%  "

$ws.Range('D20').Value = "'19.75"
$ws.Range('E20').Value = "  -8.62%  "

$ws.Range('D21').Value = "'1.007"
$ws.Range('E21').Value = "  +0.01%  "

$ws.Range('D22').Value = "'5.971"
$ws.Range('E22').Value = "  -5.56%  "

$ws.Range('D23').Value = "29.489.10"
$ws.Range('E23').Value = "  -2.65%  "

$ws.Range('D24').Value = "'11.89"
$ws.Range('E24').Value = "  -6.97%  "

$ws.Range('D25').Value = "'2.278"
$ws.Range('E25').Value = "  -2.56%  "

$ws.Range('D26').Value = "'157.62"
$ws.Range('E26').Value = "  -3.57%  "

$ws.Range('D27').Value = "'20.62"
$ws.Range('E27').Value = "  -7.34%  "

$ws.Range('D28').Value = "'6.531"
$ws.Range('E28').Value = "  -4.90%  "

$ws.Range('D29').Value = "'2.340"
$ws.Range('E29').Value = "  -8.80%  "

$ws.Range('D30').Value = "'128.17"
$ws.Range('E30').Value = "  -4.70%  "

$ws.Range('D31').Value = "'1.056"
$ws.Range('E31').Value = "  -9.26%  "

$ws.Range('D32').Value = "'0.09957"
$ws.Range('E32').Value = "  -6.12%  "

$ws.Range('D33').Value = "'1.566"
$ws.Range('E33').Value = "  -13.12%  "

$ws.Range('D34').Value = "'5.855"
$ws.Range('E34').Value = "  -7.48%  "

$ws.Range('D35').Value = "'3.787"
$ws.Range('E35').Value = "  -4.44%  "

$ws.Range('D36').Value = "'9.597"
$ws.Range('E36').Value = "  -10.89%  "

$ws.Range('D37').Value = "'0.02454"
$ws.Range('E37').Value = "  -7.51%  "

$ws.Range('D38').Value = "'0.06346"
$ws.Range('E38').Value = "  -7.81%  "

$ws.Range('D39').Value = "'1.288"
$ws.Range('E39').Value = "  -3.61%  "

$ws.Range('D40').Value = "'0.6532"
$ws.Range('E40').Value = "  -8.98%  "

$ws.Range('E41').Value = "  -8.23%  "

$ws.Range('E42').Value = "  -8.45%  "

$ws.Range('E43').Value = "  +0.08%  "

$ws.Range('D44').Value = "'0.6356"
$ws.Range('E44').Value = "  -9.01%  "

$ws.Range('D45').Value = "'2.227"
$ws.Range('E45').Value = "  -7.59%  "

$ws.Range('E46').Value = "  -9.15%  "

$ws.Range('D47').Value = "'1.271"
$ws.Range('E47').Value = "  +0.21%  "

$ws.Range('D48').Value = "'3.526"
$ws.Range('E48').Value = "  -3.14%  "

$ws.Range('D49').Value = "'0.00000000335"
$ws.Range('E49').Value = "  -3.35%  "

$ws.Range('D50').Value = "'0.06992"

$ws.Range('D51').Value = "'1.137"
$ws.Range('E51').Value = "  -6.44%  "
